# contratos-7-2011.xlsx — "fix: fixed formatting when scrapping floating point numbers"
#
# The scraper originally wrote numbers using Argentine locale formatting
# (period as thousands separator, comma as decimal separator), e.g. "10.340,00".
# This edit normalizes every such value in the "Importe" column to a plain
# floating point string with a period decimal separator, e.g. "10340.00"
# (remove thousands separators, then turn the decimal comma into a period).
# It also fixes two "Razon social" / "Nombre Fantasia" entries where a stray
# comma had been used instead of a period.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Two provider names: stray "," should be "."
# ---------------------------------------------------------------------------
$namesToFix = @(
    "FERNANDEZ MARIO H, GALLICET OSCAR M",
    "IZAGUIRRE CARLOS MARIA, MOREND MARIA ELENA Y MOREND MARIA TERESA"
)

# These values live in "Razon social" (col E) and/or "Nombre Fantasia" (col F).
foreach ($colIndex in @(5, 6)) {
    for ($r = 2; $r -le 118; $r++) {
        $cell = $ws.Cells.Item($r, $colIndex)
        $val = $cell.Value()
        if ($namesToFix -contains $val) {
            $cell.Value = $val.Replace(",", ".")
        }
    }
}

# ---------------------------------------------------------------------------
# 2) "Importe" column (H): es-AR formatted numbers -> plain decimal strings.
#    "1.234,56" -> "1234.56"  (strip "." thousands separators, then "," -> ".")
#
#    These values read like numbers, so a plain Value assignment would make
#    Excel auto-convert the cell to a numeric type. The source file keeps
#    them as text, so force a text number-format before writing, then clear
#    the format again (the value/text-type is retained) to avoid leaving a
#    stray style behind on the cell.
# ---------------------------------------------------------------------------
$importeCol = $ws.Range("H2:H118")
$importeCol.NumberFormat = "@"

for ($r = 2; $r -le 118; $r++) {
    $cell = $ws.Cells.Item($r, 8)
    $old = $cell.Value()
    $new = $old.Replace(".", "").Replace(",", ".")
    if ($new -ne $old) {
        $cell.Value = $new
    }
}

$importeCol.ClearFormats()
